$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 482
$ws.Range("C10").Value = 385.6
$ws.Range("C11").Value = 443.44
$ws.Range("C12").Value = 0
$ws.Range("C13").Value = 443.44
$ws.Range("C14").Value = 443.44
